$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add the updated EPA "What value of statistical life" hyperlink to cell B6,
# replacing the old (dead) yosemite.epa.gov link with the current epa.gov page.
$ws.Hyperlinks.Add(
    $ws.Range("B6"),
    "https://www.epa.gov/environmental-economics/mortality-risk-valuation",
    "whatvalue",
    [Type]::Missing,
    "https://www.epa.gov/environmental-economics/mortality-risk-valuation - whatvalue"
)

# The cell itself displays the full URL (with the #whatvalue fragment), which is
# distinct from the hyperlink's "display" text set above.
$ws.Range("B6").Value = "https://www.epa.gov/environmental-economics/mortality-risk-valuation#whatvalue"

# Keep the cell styled with the workbook's built-in Hyperlink style (unchanged).
$ws.Range("B6").Style = "Hyperlink"
